# Auto-generated edit script: apply Betfair odds updates for 2025-12-29 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 5.7
$ws.Range("N2").Value = 5.3
$ws.Range("O2").Value = 1.21
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 1.64
$ws.Range("R2").Value = 1.58
$ws.Range("S2").Value = 2.62
$ws.Range("T2").Value = 1.7
$ws.Range("U2").Value = 2.32
$ws.Range("W2").Value = 2.52
$ws.Range("X2").Value = 25
$ws.Range("Y2").Value = 26
$ws.Range("AA2").Value = 140
$ws.Range("AB2").Value = 11.5
$ws.Range("AC2").Value = 10.5
$ws.Range("AF2").Value = 11.5
$ws.Range("AI2").Value = 60
$ws.Range("AJ2").Value = 15.5
$ws.Range("AK2").Value = 15
$ws.Range("AL2").Value = 27
$ws.Range("AM2").Value = 80
$ws.Range("AO2").Value = 55
$ws.Range("F3").Value = 24
$ws.Range("G3").Value = 110
$ws.Range("H3").Value = 1.05
$ws.Range("I3").Value = 1.07
$ws.Range("J3").Value = 17.5
$ws.Range("L3").Value = 1.06
$ws.Range("R3").Value = 2.72
$ws.Range("S3").Value = 1.35
$ws.Range("U3").Value = 1.7
$ws.Range("X3").Value = 95
$ws.Range("AO3").Value = 29
$ws.Range("F4").Value = 2.18
$ws.Range("H4").Value = 3.9
$ws.Range("J4").Value = 2.84
$ws.Range("K4").Value = 3.25
$ws.Range("X4").Value = 14
$ws.Range("Y4").Value = 20
$ws.Range("AD4").Value = 60
$ws.Range("AF4").Value = 34
$ws.Range("AG4").Value = 23
$ws.Range("AJ4").Value = 980
$ws.Range("AK4").Value = 140
$ws.Range("G5").Value = 2.4
$ws.Range("H5").Value = 2.98
$ws.Range("I5").Value = 3.2
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 4.4
$ws.Range("Q5").Value = 1.47
$ws.Range("R5").Value = 1.74
$ws.Range("S5").Value = 2.1
$ws.Range("U5").Value = 2.82
$ws.Range("W5").Value = 1.72
$ws.Range("X5").Value = 80
$ws.Range("Y5").Value = 40
$ws.Range("Z5").Value = 85
$ws.Range("AA5").Value = 65
$ws.Range("AB5").Value = 34
$ws.Range("AI5").Value = 80
$ws.Range("AJ5").Value = 120
$ws.Range("AK5").Value = 40
$ws.Range("AL5").Value = 55
$ws.Range("AM5").Value = 580
$ws.Range("AN5").Value = 10.5
$ws.Range("G6").Value = 2.7
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 3.9
$ws.Range("J6").Value = 2.74
$ws.Range("L6").Value = 1.66
$ws.Range("AC6").Value = 7.2
$ws.Range("AE6").Value = 190
$ws.Range("I7").Value = 9
$ws.Range("J7").Value = 4.4
$ws.Range("L7").Value = 1.28
$ws.Range("N7").Value = 4.1
$ws.Range("O7").Value = 1.26
$ws.Range("P7").Value = 2.08
$ws.Range("Q7").Value = 1.76
$ws.Range("S7").Value = 2.9
$ws.Range("W7").Value = 2.96
$ws.Range("X7").Value = 20
$ws.Range("Y7").Value = 65
$ws.Range("AA7").Value = 300
$ws.Range("AB7").Value = 8.8
$ws.Range("AC7").Value = 11
$ws.Range("AD7").Value = 80
$ws.Range("AH7").Value = 46
$ws.Range("AK7").Value = 19
$ws.Range("AL7").Value = 85
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 2.66
$ws.Range("K8").Value = 3.65
$ws.Range("V8").Value = 1.53
$ws.Range("W8").Value = 1.5
$ws.Range("Y8").Value = 11.5
$ws.Range("AM8").Value = 580
$ws.Range("AN8").Value = 32
$ws.Range("F9").Value = 1.55
$ws.Range("G9").Value = 1.61
$ws.Range("H9").Value = 9.8
$ws.Range("I9").Value = 13
$ws.Range("L9").Value = 1.57
$ws.Range("N9").Value = 2.3
$ws.Range("O9").Value = 1.61
$ws.Range("P9").Value = 1.43
$ws.Range("Q9").Value = 2.84
$ws.Range("S9").Value = 6.2
$ws.Range("T9").Value = 2.84
$ws.Range("U9").Value = 1.44
$ws.Range("V9").Value = 1.08
$ws.Range("W9").Value = 2.62
$ws.Range("X9").Value = 15
$ws.Range("Y9").Value = 60
$ws.Range("AB9").Value = 8.8
$ws.Range("AC9").Value = 22
$ws.Range("AG9").Value = 36
$ws.Range("AJ9").Value = 180
$ws.Range("H10").Value = 7
$ws.Range("I10").Value = 7.2
$ws.Range("J10").Value = 3.65
$ws.Range("K10").Value = 3.7
$ws.Range("N10").Value = 2.98
$ws.Range("O10").Value = 1.48
$ws.Range("Q10").Value = 2.46
$ws.Range("R10").Value = 1.24
$ws.Range("T10").Value = 2.36
$ws.Range("X10").Value = 9.6
$ws.Range("Y10").Value = 17
$ws.Range("AF10").Value = 8
$ws.Range("AM10").Value = 240
$ws.Range("G11").Value = 1.12
$ws.Range("H11").Value = 46
$ws.Range("I11").Value = 50
$ws.Range("N11").Value = 7.6
$ws.Range("O11").Value = 1.13
$ws.Range("P11").Value = 3
$ws.Range("Q11").Value = 1.44
$ws.Range("R11").Value = 1.8
$ws.Range("S11").Value = 2.08
$ws.Range("T11").Value = 2.72
$ws.Range("U11").Value = 1.51
$ws.Range("W11").Value = 9.2
$ws.Range("Z11").Value = 550
$ws.Range("AB11").Value = 12.5
$ws.Range("AC11").Value = 30
$ws.Range("AG11").Value = 16
$ws.Range("AI11").Value = 1000
$ws.Range("AJ11").Value = 7.2
$ws.Range("AK11").Value = 16
$ws.Range("AL11").Value = 85
$ws.Range("AM11").Value = 590
$ws.Range("AN11").Value = 2.98
